# Update countries & provincias Spain
# - Banglades overtakes Colombia (rows 50/51 swap with updated Banglades figures)
# - Consejo Danes para los Refugiados overtakes Taiwan & Reunion (rows 111/112/113 shift)
# - Malasia (row 48) and Sri Lanka (row 108) figures refreshed in place
# - Header timestamp updated from 10:22 to 10:52

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 10:52"

# Row 48 - Malasia: refreshed totals, same rank
$ws.Range("B48").Value = 5780
$ws.Range("C48").Value = 38
$ws.Range("D48").Value = 3862
$ws.Range("E48").Value = 1820

# Row 50 - now Banglades (moved above Colombia) with its new totals
$ws.Range("A50").Value = "Banglades"
$ws.Range("B50").Value = 5416
$ws.Range("C50").Value = 418
$ws.Range("D50").Value = 122
$ws.Range("E50").Value = 5149
$ws.Range("F50").Value = 1
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 145

# Row 51 - now Colombia (pushed down one rank), totals unchanged
$ws.Range("A51").Value = "Colombia"
$ws.Range("B51").Value = 5142
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 1067
$ws.Range("E51").Value = 3842
$ws.Range("F51").Value = 117
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 233

# Row 108 - Sri Lanka: refreshed totals, same rank
$ws.Range("B108").Value = 467
$ws.Range("C108").Value = 15
$ws.Range("D108").Value = 120
$ws.Range("E108").Value = 340

# Row 111 - now Consejo Danes para los Refugiados (moved above Taiwan & Reunion)
$ws.Range("A111").Value = "Consejo Danes para los Refugiados"
$ws.Range("B111").Value = 442
$ws.Range("C111").Value = 26
$ws.Range("D111").Value = 50
$ws.Range("E111").Value = 364
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 28

# Row 112 - now Taiwan (pushed down one rank), totals unchanged
$ws.Range("A112").Value = "Taiwan"
$ws.Range("B112").Value = 429
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 281
$ws.Range("E112").Value = 142
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 6

# Row 113 - now Reunion (pushed down one rank), totals unchanged
$ws.Range("A113").Value = "Reunion"
$ws.Range("B113").Value = 417
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 300
$ws.Range("E113").Value = 117
$ws.Range("F113").Value = 2
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 0
